$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $val
    $ws.Range($cellRef).NumberFormat = "General"
}

$ws.Range("D2").Value = "63.229.14"
$ws.Range("E2").Value = "  +3.67%  "

$ws.Range("D3").Value = "3.457.27"
$ws.Range("E3").Value = "  +2.22%  "

Set-TextValue "D4" "0.997"
$ws.Range("E4").Value = "  -0.29%  "

Set-TextValue "D5" "584.32"
$ws.Range("E5").Value = "  +2.42%  "

Set-TextValue "D6" "148.26"
$ws.Range("E6").Value = "  +5.70%  "

$ws.Range("E7").Value = "  -0.09%  "

Set-TextValue "D8" "0.477"
$ws.Range("E8").Value = "  +0.80%  "

Set-TextValue "D9" "7.69"
$ws.Range("E9").Value = "  +0.40%  "

Set-TextValue "D10" "0.126"
$ws.Range("E10").Value = "  +3.13%  "

Set-TextValue "D11" "0.395"
$ws.Range("E11").Value = "  +3.12%  "

$ws.Range("D12").Value = "4.029.32"
$ws.Range("E12").Value = "  +1.73%  "

Set-TextValue "D13" "29.73"
$ws.Range("E13").Value = "  +7.29%  "

$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").Value = "3.464.06"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("E16").Value = "  +2.81%  "

$ws.Range("D17").Value = "63.045.86"
$ws.Range("E17").Value = "  +3.15%  "

Set-TextValue "D18" "6.26"
$ws.Range("E18").Value = "  +2.74%  "

Set-TextValue "D19" "14.29"
$ws.Range("E19").Value = "  +5.57%  "

Set-TextValue "D20" "9.38"
$ws.Range("E20").Value = "  +5.91%  "

Set-TextValue "D21" "395.52"
$ws.Range("E21").Value = "  +3.49%  "

Set-TextValue "D22" "75.40"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("E24").Value = "  -0.05%  "

Set-TextValue "D25" "0.0000120"
$ws.Range("E25").Value = "  +5.82%  "

$ws.Range("D26").Value = "3.580.47"
$ws.Range("E26").Value = "  +1.61%  "

Set-TextValue "D27" "0.192"
$ws.Range("E27").Value = "  -0.17%  "

Set-TextValue "D28" "7.79"
$ws.Range("E28").Value = "  +8.69%  "

$ws.Range("E29").Value = "  -0.11%  "

Set-TextValue "D30" "8.19"
$ws.Range("E30").Value = "  +3.22%  "

$ws.Range("E31").Value = "  +1.47%  "

Set-TextValue "D32" "1.42"
$ws.Range("E32").Value = "  +5.06%  "

Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  -0.04%  "

Set-TextValue "D34" "23.83"
$ws.Range("E34").Value = "  +2.78%  "

Set-TextValue "D35" "5.34"
$ws.Range("E35").Value = "  +8.16%  "

Set-TextValue "D36" "7.15"
$ws.Range("E36").Value = "  +3.69%  "

$ws.Range("E37").Value = "  +8.82%  "

Set-TextValue "D38" "169.70"
$ws.Range("E38").Value = "  +1.93%  "

Set-TextValue "D39" "31.25"
$ws.Range("E39").Value = "  +19.99%  "

$ws.Range("D40").Value = "3.477.41"
$ws.Range("E40").Value = "  +1.61%  "

Set-TextValue "D41" "0.0772"
$ws.Range("E41").Value = "  +1.17%  "

Set-TextValue "D42" "0.794"
$ws.Range("E42").Value = "  +2.01%  "

Set-TextValue "D43" "42.96"
$ws.Range("E43").Value = "  +0.96%  "

$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D44" "1.75"
$ws.Range("E44").Value = "  +7.03%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D45" "4.49"
$ws.Range("E45").Value = "  +3.39%  "

Set-TextValue "D46" "1.21"
$ws.Range("E46").Value = "  +8.67%  "

$ws.Range("D47").Value = "2.569.19"
$ws.Range("E47").Value = "  +5.04%  "

Set-TextValue "D48" "23.66"
$ws.Range("E48").Value = "  +4.25%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "2.23"
$ws.Range("E49").Value = "  +7.99%  "

$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D50" "6.75"
$ws.Range("E50").Value = "  +2.14%  "

$ws.Range("E51").Value = "  -0.05%  "
